$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.04790737368050202
$ws.Range("D2").Value = 0.08984926657176118
$ws.Range("E2").Value = 0.1134667721616456
$ws.Range("F2").Value = 2.191584079143524
$ws.Range("G2").Value = 1.547441977182544
$ws.Range("H2").Value = 1.383694366216645
$ws.Range("J2").Value = 0.1593564290105576
$ws.Range("K2").Value = 1.574579318661733
$ws.Range("M2").Value = 0.4711100814171303

$ws.Range("C3").Value = 0.04260070930591553
$ws.Range("D3").Value = 0.08884552614073726
$ws.Range("E3").Value = 0.113453344488164
$ws.Range("F3").Value = 2.192549122902008
$ws.Range("G3").Value = 1.545299445426807
$ws.Range("H3").Value = 1.390069314688546
$ws.Range("J3").Value = 0.1602241914952955
$ws.Range("K3").Value = 1.437292594492931
$ws.Range("M3").Value = 0.4451127177002476

$ws.Range("C4").Value = 0.03935978376252081
$ws.Range("D4").Value = 0.08824260691808661
$ws.Range("E4").Value = 0.1134868023295592
$ws.Range("F4").Value = 2.194599706022132
$ws.Range("G4").Value = 1.54520239302569
$ws.Range("H4").Value = 1.394810783348376
$ws.Range("J4").Value = 0.1608353916333947
$ws.Range("K4").Value = 1.35340631648873
$ws.Range("M4").Value = 0.4293336202403708

$ws.Range("C5").Value = 0.03804337332059049
$ws.Range("D5").Value = 0.08800031083804782
$ws.Range("E5").Value = 0.113510935269856
$ws.Range("F5").Value = 2.195800942765771
$ws.Range("G5").Value = 1.54546794732579
$ws.Range("H5").Value = 1.396950535201128
$ws.Range("J5").Value = 0.1611041528057306
$ws.Range("K5").Value = 1.319325205370802
$ws.Range("M5").Value = 0.4229497474410522

$ws.Range("C6").Value = 0.03782504135922693
$ws.Range("D6").Value = 0.08796028401194178
$ws.Range("E6").Value = 0.1135155769791467
$ws.Range("F6").Value = 2.196022458971086
$ws.Range("G6").Value = 1.545530432230407
$ws.Range("H6").Value = 1.397318362720355
$ws.Range("J6").Value = 0.1611499691776572
$ws.Range("K6").Value = 1.313672314742689
$ws.Range("M6").Value = 0.4218925066646548

$ws.Range("C7").Value = 0.03934201291181694
$ws.Range("D7").Value = 0.08823932542253488
$ws.Range("E7").Value = 0.113487085269858
$ws.Range("F7").Value = 2.194614427367569
$ws.Range("G7").Value = 1.545204740748744
$ws.Range("H7").Value = 1.394838800973147
$ws.Range("J7").Value = 0.1608389365340912
$ws.Range("K7").Value = 1.352946267975966
$ws.Range("M7").Value = 0.429247337716383

$ws.Range("C8").Value = 0.0460739774795087
$ws.Range("D8").Value = 0.08950041906479811
$ws.Range("E8").Value = 0.1134534939473593
$ws.Range("F8").Value = 2.191613583281153
$ws.Range("G8").Value = 1.546449556215805
$ws.Range("H8").Value = 1.385720486308031
$ws.Range("J8").Value = 0.1596393514985515
$ws.Range("K8").Value = 1.527158382938751
$ws.Range("M8").Value = 0.4621081904478714

$ws.Range("C9").Value = 0.05941756285893973
$ws.Range("D9").Value = 0.09207830819200069
$ws.Range("E9").Value = 0.1137180920129985
$ws.Range("F9").Value = 2.197347209229918
$ws.Range("G9").Value = 1.558619333340232
$ws.Range("H9").Value = 1.374424382873173
$ws.Range("J9").Value = 0.1579099360062131
$ws.Range("K9").Value = 1.87202363694621
$ws.Range("M9").Value = 0.5280020322065582

$ws.Range("C10").Value = 0.06931467748201214
$ws.Range("D10").Value = 0.09403476808722644
$ws.Range("E10").Value = 0.1141135617527063
$ws.Range("F10").Value = 2.208715115290317
$ws.Range("G10").Value = 1.573578929368296
$ws.Range("H10").Value = 1.370170370931476
$ws.Range("J10").Value = 0.157020581008176
$ws.Range("K10").Value = 2.12739321531086
$ws.Range("M10").Value = 0.5773048456172774

$ws.Range("C11").Value = 0.07383900106879082
$ws.Range("D11").Value = 0.09493809946535947
$ws.Range("E11").Value = 0.1143370643729895
$ws.Range("F11").Value = 2.215456405632665
$ws.Range("G11").Value = 1.581710445348307
$ws.Range("H11").Value = 1.369120314797499
$ws.Range("J11").Value = 0.1566991049330753
$ws.Range("K11").Value = 2.244008844402288
$ws.Range("M11").Value = 0.5999288002571035

$ws.Range("C12").Value = 0.07555552423575307
$ws.Range("D12").Value = 0.09528205676339496
$ws.Range("E12").Value = 0.1144279608499374
$ws.Range("F12").Value = 2.218236153033828
$ws.Range("G12").Value = 1.584981859311625
$ws.Range("H12").Value = 1.368850503918424
$ws.Range("J12").Value = 0.1565893455539609
$ws.Range("K12").Value = 2.288232469414936
$ws.Range("M12").Value = 0.6085240741815738

$ws.Range("C13").Value = 0.07518569389962693
$ws.Range("D13").Value = 0.09520789604491853
$ws.Range("E13").Value = 0.1144081063521831
$ws.Range("F13").Value = 2.217627370968941
$ws.Range("G13").Value = 1.584268730929494
$ws.Range("H13").Value = 1.368902918861778
$ws.Range("J13").Value = 0.1566124510940057
$ws.Range("K13").Value = 2.278705288473247
$ws.Range("M13").Value = 0.6066716809374384

$ws.Range("C14").Value = 0.0739801546350094
$ws.Range("D14").Value = 0.09496635943065712
$ws.Range("E14").Value = 0.1143444170535837
$ws.Range("F14").Value = 2.215680541426693
$ws.Range("G14").Value = 1.581975726465032
$ws.Range("H14").Value = 1.369095552407941
$ws.Range("J14").Value = 0.156689834774653
$ws.Range("K14").Value = 2.247645871917825
$ws.Range("M14").Value = 0.6006353764739032

$ws.Range("C15").Value = 0.0732421541269872
$ws.Range("D15").Value = 0.0948186558891706
$ws.Range("E15").Value = 0.1143062206107608
$ws.Range("F15").Value = 2.214517645877891
$ws.Range("G15").Value = 1.580596266732584
$ws.Range("H15").Value = 1.369230208509805
$ws.Range("J15").Value = 0.156738795047076
$ws.Range("K15").Value = 2.228629381559358
$ws.Range("M15").Value = 0.5969416191233989

$ws.Range("C16").Value = 0.06901945381554242
$ws.Range("D16").Value = 0.09397599876248819
$ws.Range("E16").Value = 0.1140998319316608
$ws.Range("F16").Value = 2.208306250969969
$ws.Range("G16").Value = 1.573074329462827
$ws.Range("H16").Value = 1.37025685176539
$ws.Range("J16").Value = 0.1570432646756217
$ws.Range("K16").Value = 2.119781061957497
$ws.Range("M16").Value = 0.5758302501718333

$ws.Range("C17").Value = 0.06643467878761555
$ws.Range("D17").Value = 0.09346244706022588
$ws.Range("E17").Value = 0.1139843800768396
$ws.Range("F17").Value = 2.204898725351072
$ws.Range("G17").Value = 1.568800605254864
$ws.Range("H17").Value = 1.37111374919678
$ws.Range("J17").Value = 0.15725134928784
$ws.Range("K17").Value = 2.053120247322227
$ws.Range("M17").Value = 0.5629292055953812

$ws.Range("C18").Value = 0.06495005733962955
$ws.Range("D18").Value = 0.09316832201454872
$ws.Range("E18").Value = 0.1139220797160441
$ws.Range("F18").Value = 2.203086507448688
$ws.Range("G18").Value = 1.56646722393981
$ws.Range("H18").Value = 1.371689892416811
$ws.Range("J18").Value = 0.1573788537761871
$ws.Range("K18").Value = 2.014820780400214
$ws.Range("M18").Value = 0.555527305893996

$ws.Range("C19").Value = 0.06444774344080884
$ws.Range("D19").Value = 0.09306895306406915
$ws.Range("E19").Value = 0.1139016911121971
$ws.Range("F19").Value = 2.202498252677685
$ws.Range("G19").Value = 1.565698562910768
$ws.Range("H19").Value = 1.371899251248152
$ws.Range("J19").Value = 0.1574233668155074
$ws.Range("K19").Value = 2.001860478190451
$ws.Range("M19").Value = 0.553024319818924

$ws.Range("C20").Value = 0.06670961735819958
$ws.Range("D20").Value = 0.09351698572018563
$ws.Range("E20").Value = 0.1139962453904602
$ws.Range("F20").Value = 2.205246167133581
$ws.Range("G20").Value = 1.569242628617246
$ws.Range("H20").Value = 1.371013908099144
$ws.Range("J20").Value = 0.1572283888173729
$ws.Range("K20").Value = 2.060212053877251
$ws.Range("M20").Value = 0.5643006360371743

$ws.Range("C21").Value = 0.07433416157816453
$ws.Range("D21").Value = 0.09503725366751326
$ws.Range("E21").Value = 0.1143629542956504
$ws.Range("F21").Value = 2.216246202976464
$ws.Range("G21").Value = 1.582644009917914
$ws.Range("H21").Value = 1.369035497839945
$ws.Range("J21").Value = 0.1566667800480701
$ws.Range("K21").Value = 2.256767046617995
$ws.Range("M21").Value = 0.6024076241957772

$ws.Range("C22").Value = 0.07933628413285021
$ws.Range("D22").Value = 0.09604181171582837
$ws.Range("E22").Value = 0.1146391090550658
$ws.Range("F22").Value = 2.224758865722436
$ws.Range("G22").Value = 1.592523474045322
$ws.Range("H22").Value = 1.368487795797876
$ws.Range("J22").Value = 0.1563695582603231
$ws.Range("K22").Value = 2.385598907170333
$ws.Range("M22").Value = 0.6274763406284904

$ws.Range("C23").Value = 0.07666478636791396
$ws.Range("D23").Value = 0.09550466608924069
$ws.Range("E23").Value = 0.1144883841691851
$ws.Range("F23").Value = 2.22009399591181
$ws.Range("G23").Value = 1.587147571990045
$ws.Range("H23").Value = 1.368711738471887
$ws.Range("J23").Value = 0.1565217932172942
$ws.Range("K23").Value = 2.316805026001475
$ws.Range("M23").Value = 0.6140817534885201

$ws.Range("C24").Value = 0.06658531327713035
$ws.Range("D24").Value = 0.09349232529261542
$ws.Range("E24").Value = 0.113990868390399
$ws.Range("F24").Value = 2.205088631480805
$ws.Range("G24").Value = 1.569042404840076
$ws.Range("H24").Value = 1.371058786196187
$ws.Range("J24").Value = 0.1572387447223846
$ws.Range("K24").Value = 2.057005771098432
$ws.Range("M24").Value = 0.5636805653485482

$ws.Range("C25").Value = 0.05579183542131716
$ws.Range("D25").Value = 0.09136984539028958
$ws.Range("E25").Value = 0.1136111739006189
$ws.Range("F25").Value = 2.194544815686399
$ws.Range("G25").Value = 1.554276529663497
$ws.Range("H25").Value = 1.37677213917658
$ws.Range("J25").Value = 0.1583109450440112
$ws.Range("K25").Value = 1.778379666027092
$ws.Range("M25").Value = 0.5100199467522941
